# Mark four existing wishlist items as reserved/bought ("Y" in column E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = "Y"
$ws.Range("E22").Value = "Y"
$ws.Range("E31").Value = "Y"
$ws.Range("E42").Value = "Y"

# Append three new wishlist rows (45-47) with Name / Image / Link / Price.
# Values are written in the same per-row column order the source data was
# pasted in so the workbook's shared-string table matches the upload.

# Row 45: The Leopard (Criterion Blu-ray)
$ws.Range("C45").Value = "https://www.criterion.com/films/790-the-leopard"
$ws.Range("D45").Value = "40 USD"
$ws.Range("B45").Value = "https://s3.amazonaws.com/criterion-production/films/96833308352eedd79b50b2de8016a6c5/7IQOqWTrRuJg90rWrJZbK1VXI7aDcR_large.jpg"
$ws.Range("A45").Value = "The Leopard"

# Row 46: Heat: Pedal to the Metal (English Edition)
$ws.Range("A46").Value = "Heat: Pedal to the Metal (English Edition)"
$ws.Range("B46").Value = "https://regatuljocurilor.ro/74113-large_default/heat-pedal-to-the-metal.jpg"
$ws.Range("C46").Value = "https://regatuljocurilor.ro/ro/acasa/heat-pedal-to-the-metal"
$ws.Range("D46").Value = "399 RON"

# Row 47: Sky Team (2024 Romanian Edition)
$ws.Range("B47").Value = "https://regatuljocurilor.ro/134548-large_default/sky-team-romanian-edition.jpg"
$ws.Range("C47").Value = "https://regatuljocurilor.ro/ro/acasa/sky-team-romanian-edition"
$ws.Range("D47").Value = "169 RON"
$ws.Range("A47").Value = "Sky Team (2024 Romanian Edition)"

# Match the saved file's final selection (last edited cell).
$null = $ws.Range("A47").Select()
